$d = $word.ActiveDocument

# 1) "Yo, " + "FIRST_" + "NAME LAST" + "_" + "NAME , identificado..." is collapsed
#    into a single run with the same (unchanged) concatenated text. Doing a
#    Find/Replace across the whole phrase with identical text merges the
#    previously split runs into one, matching the target XML.
$d.Content.Find.Execute(
    "Yo, FIRST_NAME LAST_NAME , identificado con c",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Yo, FIRST_NAME LAST_NAME , identificado con c", 2) | Out-Null

# 2) Fill in more of the closing date sentence: the placeholder "DAY" becomes
#    "DATE", and the "del mes de  MONTH de YEAR" wording is tightened to
#    "del mes MONTH del YEAR" (drop the stray double space / extra "de", add
#    the missing "l" so "de YEAR" reads "del YEAR").
$d.Content.Find.Execute(
    "DAY", $true, $false, $false, $false, $false, $true, 1, $false,
    "DATE", 2) | Out-Null

$d.Content.Find.Execute(
    " del mes de  MONTH", $true, $false, $false, $false, $false, $true, 1, $false,
    " del mes MONTH", 2) | Out-Null

$d.Content.Find.Execute(
    "MONTH de YEAR", $true, $false, $false, $false, $false, $true, 1, $false,
    "MONTH del YEAR", 2) | Out-Null

# 3) "Nombre: NAME LASTNAME" -> "Nombre: FIRST_NAME LAST_NAME"
$d.Content.Find.Execute(
    "Nombre: NAME LASTNAME", $true, $false, $false, $false, $false, $true, 1, $false,
    "Nombre: FIRST_NAME LAST_NAME", 2) | Out-Null
